$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.870.79"
Set-TextValue "E2" "  +0.54%  "
Set-TextValue "D3" "1.643.01"
Set-TextValue "E3" "  +0.59%  "
Set-TextValue "E4" "  -0.80%  "
Set-TextValue "D5" "216.69"
Set-TextValue "E5" "  -0.61%  "
Set-TextValue "D6" "0.506"
Set-TextValue "E6" "  +1.70%  "
Set-TextValue "E7" "  -0.65%  "
Set-TextValue "E8" "  +1.61%  "
Set-TextValue "E9" "  +0.43%  "
Set-TextValue "D10" "19.78"
Set-TextValue "E10" "  +4.38%  "
Set-TextValue "D11" "0.0847"
Set-TextValue "E11" "  +0.52%  "
Set-TextValue "D12" "1.873.58"
Set-TextValue "E12" "  +0.68%  "
Set-TextValue "D13" "1.624.89"
Set-TextValue "E13" "  -0.30%  "
Set-TextValue "D14" "4.12"
Set-TextValue "E14" "  +0.26%  "
Set-TextValue "E15" "  +1.24%  "
Set-TextValue "D16" "66.17"
Set-TextValue "E16" "  +3.38%  "
Set-TextValue "D17" "26.901.20"
Set-TextValue "E17" "  +0.80%  "
Set-TextValue "E18" "  +0.61%  "
Set-TextValue "D19" "219.54"
Set-TextValue "E19" "  +4.04%  "
Set-TextValue "E20" "  -0.78%  "
Set-TextValue "D21" "4.37"
Set-TextValue "E21" "  +1.38%  "
Set-TextValue "D22" "6.61"
Set-TextValue "E22" "  +6.91%  "
Set-TextValue "E23" "  +3.42%  "
Set-TextValue "D24" "9.17"
Set-TextValue "E24" "  -0.16%  "
Set-TextValue "D25" "145.92"
Set-TextValue "E25" "  -0.68%  "
Set-TextValue "E26" "  -0.78%  "
Set-TextValue "D27" "7.41"
Set-TextValue "E27" "  +6.06%  "
Set-TextValue "E28" "  +1.50%  "
Set-TextValue "D29" "15.81"
Set-TextValue "E29" "  +1.66%  "
Set-TextValue "E30" "  +0.77%  "
Set-TextValue "E31" "  -0.53%  "
Set-TextValue "D32" "3.34"
Set-TextValue "E32" "  -0.71%  "
Set-TextValue "E33" "  +1.81%  "
Set-TextValue "E34" "  +2.47%  "
Set-TextValue "B35" "HuobiToken"
Set-TextValue "C35" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D35" "2.44"
Set-TextValue "E35" "  -0.29%  "
Set-TextValue "B36" "Maker"
Set-TextValue "C36" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D36" "1.247.13"
Set-TextValue "E36" "  -1.26%  "
Set-TextValue "E37" "  +1.30%  "
Set-TextValue "D39" "0.831"
Set-TextValue "E39" "  +3.54%  "
Set-TextValue "E40" "  -0.78%  "
Set-TextValue "E41" "  +0.94%  "
Set-TextValue "D42" "5.35"
Set-TextValue "E42" "  +1.91%  "
Set-TextValue "D43" "1.784.93"
Set-TextValue "E43" "  +0.75%  "
Set-TextValue "E44" "  -2.58%  "
Set-TextValue "E45" "  +1.83%  "
Set-TextValue "D46" "91.51"
Set-TextValue "E46" "  +0.17%  "
Set-TextValue "E47" "  +0.85%  "
Set-TextValue "E48" "  +16.86%  "
Set-TextValue "E49" "  -0.07%  "
Set-TextValue "D50" "0.0972"
Set-TextValue "E50" "  +1.80%  "
Set-TextValue "D51" "7.61"
Set-TextValue "E51" "  +1.91%  "
